# metro_budget_exercise.xlsx edit
# "Question 6, dropdown B87 done, beginning to fill in budget table"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# ---------------------------------------------------------------
# Question headers: retitle + reorder shared strings so the final
# table matches "Question N (<method>)" naming. We touch A72/A63/A54
# in that order so the newly-appended shared strings land in the
# sequence: Question 5 (Index + Match), Question 4 (XLOOKUP),
# Question 3 (VLOOKUP).
# ---------------------------------------------------------------
$ws.Range("A72").Value = "Question 5 (Index + Match)"
$ws.Range("A63").Value = "Question 4 (XLOOKUP)"
$ws.Range("A54").Value = "Question 3 (VLOOKUP)"

# ---------------------------------------------------------------
# Question 3 (VLOOKUP) table, rows 56-61: convert the C column from
# XLOOKUP to VLOOKUP and add the new D column (FY19_diff lookups).
# ---------------------------------------------------------------
$ws.Range("C56").Formula = "=VLOOKUP(A10,A2:I52,9)"
$ws.Range("D56").Formula = "=VLOOKUP(A10,A2:N52,14)"

$ws.Range("C57").Formula = "=VLOOKUP(A11,A2:I52,9)"
$ws.Range("D57").Formula = "=VLOOKUP(A11,A2:N52,14)"

$ws.Range("C58").Formula = "=VLOOKUP(A18,A2:I52,9)"
$ws.Range("D58").Formula = "=VLOOKUP(A18,A2:N52,14)"

$ws.Range("C59").Formula = "=VLOOKUP(A24,A2:I52,9)"
$ws.Range("D59").Formula = "=VLOOKUP(A24,A2:N52,14)"

$ws.Range("C60").Formula = "=VLOOKUP(A25,A2:I52,9)"
$ws.Range("D60").Formula = "=VLOOKUP(A25,A2:N52,14)"

$ws.Range("C61").Formula = "=VLOOKUP(A41,A2:I52,9)"
$ws.Range("D61").Formula = "=VLOOKUP(A41,A2:N52,14)"

# ---------------------------------------------------------------
# Question 4 (XLOOKUP) table, rows 65-70: fill in B/C/D columns
# (FY17_diff, FY18_diff, FY19_diff) for each department.
# ---------------------------------------------------------------
$ws.Range("B65").Formula = '=_xlfn.XLOOKUP("Community Education Commission", A2:A52, D2:D52)'
$ws.Range("C65").Formula = '=_xlfn.XLOOKUP("Community Education Commission", A2:A52, I2:I52)'
$ws.Range("D65").Formula = '=_xlfn.XLOOKUP("Community Education Commission", A2:A52, N2:N52)'

$ws.Range("B66").Formula = '=_xlfn.XLOOKUP("Community Oversight Board", A2:A52, D2:D52)'
$ws.Range("C66").Formula = '=_xlfn.XLOOKUP("Community Oversight Board", A2:A52, I2:I52)'
$ws.Range("D66").Formula = '=_xlfn.XLOOKUP("Community Oversight Board", A2:A52, N2:N52)'

$ws.Range("B67").Formula = '=_xlfn.XLOOKUP("Election Commission", A2:A52, D2:D52)'
$ws.Range("C67").Formula = '=_xlfn.XLOOKUP("Election Commission", A2:A52, I2:I52)'
$ws.Range("D67").Formula = '=_xlfn.XLOOKUP("Election Commission", A2:A52, N2:N52)'

$ws.Range("B68").Formula = '=_xlfn.XLOOKUP("Historical Commission", A2:A52, D2:D52)'
$ws.Range("C68").Formula = '=_xlfn.XLOOKUP("Historical Commission", A2:A52, I2:I52)'
$ws.Range("D68").Formula = '=_xlfn.XLOOKUP("Historical Commission", A2:A52, N2:N52)'

$ws.Range("B69").Formula = '=_xlfn.XLOOKUP("Human Relations Commission", A2:A52, D2:D52)'
$ws.Range("C69").Formula = '=_xlfn.XLOOKUP("Human Relations Commission", A2:A52, I2:I52)'
$ws.Range("D69").Formula = '=_xlfn.XLOOKUP("Human Relations Commission", A2:A52, N2:N52)'

$ws.Range("B70").Formula = '=_xlfn.XLOOKUP("Planning Commission", A2:A52, D2:D52)'
$ws.Range("C70").Formula = '=_xlfn.XLOOKUP("Planning Commission", A2:A52, I2:I52)'
$ws.Range("D70").Formula = '=_xlfn.XLOOKUP("Planning Commission", A2:A52, N2:N52)'

# ---------------------------------------------------------------
# Question 5 (Index + Match) table, rows 74-79: same layout using
# INDEX/MATCH instead of XLOOKUP.
# ---------------------------------------------------------------
$ws.Range("B74").Formula = '=INDEX(D2:D52,MATCH("Community Education Commission",A2:A52,0))'
$ws.Range("C74").Formula = '=INDEX(I2:I52,MATCH("Community Education Commission",A2:A52,0))'
$ws.Range("D74").Formula = '=INDEX(N2:N52,MATCH("Community Education Commission",A2:A52,0))'

$ws.Range("B75").Formula = '=INDEX(D2:D52,MATCH("Community Oversight Board",A2:A52,0))'
$ws.Range("C75").Formula = '=INDEX(I2:I52,MATCH("Community Oversight Board",A2:A52,0))'
$ws.Range("D75").Formula = '=INDEX(N2:N52,MATCH("Community Oversight Board",A2:A52,0))'

$ws.Range("B76").Formula = '=INDEX(D2:D52,MATCH("Election Commission",A2:A52,0))'
$ws.Range("C76").Formula = '=INDEX(I2:I52,MATCH("Election Commission",A2:A52,0))'
$ws.Range("D76").Formula = '=INDEX(N2:N52,MATCH("Election Commission",A2:A52,0))'

$ws.Range("B77").Formula = '=INDEX(D2:D52,MATCH("Historical Commission",A2:A52,0))'
$ws.Range("C77").Formula = '=INDEX(I2:I52,MATCH("Historical Commission",A2:A52,0))'
$ws.Range("D77").Formula = '=INDEX(N2:N52,MATCH("Historical Commission",A2:A52,0))'

$ws.Range("B78").Formula = '=INDEX(D2:D52,MATCH("Human Relations Commission",A2:A52,0))'
$ws.Range("C78").Formula = '=INDEX(I2:I52,MATCH("Human Relations Commission",A2:A52,0))'
$ws.Range("D78").Formula = '=INDEX(N2:N52,MATCH("Human Relations Commission",A2:A52,0))'

$ws.Range("B79").Formula = '=INDEX(D2:D52,MATCH("Planning Commission",A2:A52,0))'
$ws.Range("C79").Formula = '=INDEX(I2:I52,MATCH("Planning Commission",A2:A52,0))'
$ws.Range("D79").Formula = '=INDEX(N2:N52,MATCH("Planning Commission",A2:A52,0))'

# ---------------------------------------------------------------
# Question 6 budget table: start filling it in with the first
# INDEX/MATCH lookup (Administrative budget).
# ---------------------------------------------------------------
$ws.Range("B84").Formula = '=INDEX(B2:B52,MATCH("Administrative",A2:A52,0))'

# Dropdown list validation on B87 (same list as A83's department picker).
$ws.Range("B87").Validation.Add(3, 1, 1, "=`$A`$2:`$A`$52")

# ---------------------------------------------------------------
# Selection / scroll position bookkeeping to mirror the saved view.
# ---------------------------------------------------------------
$ws.Range("A55").Select()
